# Trade #54 closed at 2026-02-17 12:49:35 - unknown UNKNOWN +0.000%
#
# Updates the running totals on the "Summary" and "Strategy Status" sheets,
# and appends the newly-closed trade (#54) as a new row to both the
# "All Trades" log and the per-strategy "MarketMaking" log.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: refresh the aggregate stats after the new trade closed
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.19   # Current Capital
$summary.Range("B4").Value = 0.18      # Total P&L $
$summary.Range("B5").Value = 0.07      # Total P&L %
$summary.Range("B6").Value = 54        # Total Trades
$summary.Range("B7").Value = 23        # Winning Trades
$summary.Range("B9").Value = 42.59     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet: row 4 holds the "MarketMaking" strategy totals
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.19   # Capital
$status.Range("D4").Value = 54       # Trades
$status.Range("E4").Value = 0.18     # P&L $
$status.Range("F4").Value = 0.19     # P&L %
$status.Range("G4").Value = 42.59    # Win Rate %

# ---------------------------------------------------------------------
# Append the new trade (row 55) to a trade-log sheet. Column B holds a
# plain "YYYY-MM-DD" text date; writing it with a leading apostrophe
# forces Excel to keep it as literal text instead of auto-converting it
# to a date serial number (matching how the existing rows are stored).
# ---------------------------------------------------------------------
function Add-Trade55($ws) {
    $ws.Cells.Item(55, 1).Value = 54
    $ws.Cells.Item(55, 2).Value = "'2026-02-17"
    $ws.Cells.Item(55, 3).Value = "12:49:28"
    $ws.Cells.Item(55, 4).Value = "MarketMaking"
    $ws.Cells.Item(55, 5).Value = "UP"
    $ws.Cells.Item(55, 6).Value = 0.98
    $ws.Cells.Item(55, 7).Value = 0.99
    $ws.Cells.Item(55, 8).Value = "CLOSED"
    $ws.Cells.Item(55, 9).Value = 1.0204
    $ws.Cells.Item(55, 10).Value = 0.01
    $ws.Cells.Item(55, 11).Value = 100.19
    $ws.Cells.Item(55, 12).Value = 0
    $ws.Cells.Item(55, 13).Value = 0
    $ws.Cells.Item(55, 14).Value = 0.6
    $ws.Cells.Item(55, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(55, 16).Value = "early_exit"
    $ws.Cells.Item(55, 17).Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade55 $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade55 $marketMaking
